$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Protocol")

# Update the GUI-mockup progress note for Stefan Herbst (3. Sprint, first half -> column E)
$ws.Range("E3").Value = "Fertigstellung GUI Mockups; Start Erstellung der GUI-Struktur in WPF"

# Add the new progress note for the second half of 3. Sprint (column F)
$ws.Range("F3").Value = "Fertigstellung GUI-Struktur in WPF; Start Implementierung der Business Logic"

# Update progress percentages
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.2

# Reflect the new selection cell as last edited in Excel
$null = $ws.Range("F6").Select()
